$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '37.062.49'
$ws.Range("E2").Value = '  -1.30%  '
$ws.Range("D3").Value = '2.015.68'
$ws.Range("E3").Value = '  -2.68%  '
$ws.Range("D4").Value = "'1.00"
$ws.Range("D4").Style = "Normal"
$ws.Range("E4").Value = '  -0.12%  '
$ws.Range("D5").Value = "'226.01"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = '  -2.45%  '
$ws.Range("E6").Value = '  -2.60%  '
$ws.Range("E7").Value = '  +0.04%  '
$ws.Range("D8").Value = "'54.87"
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = '  -5.29%  '
$ws.Range("D9").Value = "'0.377"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = '  -2.78%  '
$ws.Range("E10").Value = '  +0.88%  '
$ws.Range("E11").Value = '  -5.05%  '
$ws.Range("D12").Value = '2.313.13'
$ws.Range("D13").Value = "'14.14"
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = '  -4.12%  '
$ws.Range("D14").Value = "'20.23"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = '  -4.52%  '
$ws.Range("E15").Value = '  -3.08%  '
$ws.Range("D16").Value = "'5.12"
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = '  -3.77%  '
$ws.Range("D17").Value = '2.009.52'
$ws.Range("E17").Value = '  -3.25%  '
$ws.Range("D18").Value = '37.013.87'
$ws.Range("E18").Value = '  -1.32%  '
$ws.Range("D19").Value = "'6.17"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = '  +0.28%  '
$ws.Range("D20").Value = "'68.78"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = '  -1.69%  '
$ws.Range("D21").Value = '0.0₃0815'
$ws.Range("E21").Value = '  -1.29%  '
$ws.Range("D22").Value = "'223.01"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = '  -1.80%  '
$ws.Range("E23").Value = '  -0.04%  '
$ws.Range("E24").Value = '  +1.52%  '
$ws.Range("E25").Value = '  -6.14%  '
$ws.Range("D26").Value = "'165.93"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = '  -2.11%  '
$ws.Range("D27").Value = "'9.16"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = '  -7.14%  '
$ws.Range("E28").Value = '  -0.11%  '
$ws.Range("D29").Value = "'18.69"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = '  -3.41%  '
$ws.Range("E30").Value = '  -4.87%  '
$ws.Range("D31").Value = "'0.117"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = '  -3.89%  '
$ws.Range("D32").Value = "'4.52"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = '  -0.82%  '
$ws.Range("D33").Value = "'0.0612"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = '  -2.29%  '
$ws.Range("E34").Value = '  -4.75%  '
$ws.Range("D35").Value = "'2.35"
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = '  -7.04%  '
$ws.Range("E36").Value = '  +0.62%  '
$ws.Range("E37").Value = '  -0.22%  '
$ws.Range("E38").Value = '  -4.24%  '
$ws.Range("D39").Value = "'5.29"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = '  -0.44%  '
$ws.Range("D40").Value = '1.481.73'
$ws.Range("D41").Value = "'0.0216"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = '  -4.84%  '
$ws.Range("D42").Value = "'94.85"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = '  -3.37%  '
$ws.Range("D43").Value = "'0.0914"
$ws.Range("D43").Style = "Normal"
$ws.Range("D44").Value = "'16.32"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = '  -1.50%  '
$ws.Range("E45").Value = '  -5.05%  '
$ws.Range("E46").Value = '  -5.24%  '
$ws.Range("D47").Value = "'7.20"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = '  -0.68%  '
$ws.Range("E48").Value = '  -2.59%  '
$ws.Range("D49").Value = "'2.92"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = '  -1.24%  '
$ws.Range("D50").Value = '2.198.61'
$ws.Range("E50").Value = '  -2.76%  '
$ws.Range("B51").Value = 'MultiversX'
$ws.Range("C51").Value = 'https://coinranking.com/coin/omwkOTglq+multiversx-egld'
$ws.Range("D51").Value = "'44.39"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = '  -2.79%  '
